$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update the "Estatus" (column F) values for the affected tasks.
$ws.Range("F6").Value = "Hecho"
$ws.Range("F7").Value = "Hecho"
$ws.Range("F8").Value = "En proceso"
$ws.Range("F9").Value = "Hecho"
$ws.Range("F11").Value = "Hecho"

# Row 11: 4 hours consumed on day 1 (K11), which ripples through the
# remaining day formulas (L11, O11, R11, U11, X11, AA11, AD11, AG11, AJ11,
# AM11, AP11, AS11, AV11, AY11, AZ11, BA11).
$ws.Range("K11").Value = 4

# Move the active selection on the frozen (bottom-right) pane to F8.
$ws.Activate()
$ws.Range("F8").Select()
